# Updates the cryptos list (Price / Volume(1h) columns) to the latest
# scraped values, as produced by the scheduled GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column cells are stored as text in the workbook (thousands are
# separated with '.', e.g. "70.187.02"), not as numbers. Without the
# leading apostrophe, Excel's type-inference would read plain decimals
# (e.g. "53.90", "12.30") as a Number and silently drop the significant
# trailing zero, so always enter these as text - same as typing
# '53.90 into the grid by hand.
function Set-Price($row, $value) {
    $ws.Cells.Item($row, 4).Value = "'" + $value
}

function Set-Volume($row, $value) {
    $ws.Cells.Item($row, 5).Value = "  $value  "
}

Set-Price 2 "70.187.02"

Set-Price 3 "3.608.90"
Set-Volume 3 "+2.84%"

Set-Volume 4 "+0.14%"

Set-Price 5 "603.15"
Set-Volume 5 "+0.43%"

Set-Price 6 "196.33"
Set-Volume 6 "-0.04%"

Set-Volume 7 "+0.37%"

Set-Volume 8 "+0.06%"

Set-Price 9 "0.207"
Set-Volume 9 "-1.23%"

Set-Price 10 "0.649"
Set-Volume 10 "-0.71%"

Set-Price 11 "53.90"
Set-Volume 11 "-0.35%"

Set-Volume 12 "+1.23%"

Set-Volume 13 "+0.33%"

Set-Price 14 "4.174.29"
Set-Volume 14 "+2.83%"

Set-Price 15 "13.18"
Set-Volume 15 "+4.57%"

Set-Price 16 "592.08"

Set-Price 17 "70.281.97"
Set-Volume 17 "+0.35%"

Set-Price 18 "19.18"
Set-Volume 18 "+1.00%"

Set-Price 19 "3.607.47"
Set-Volume 19 "+3.03%"

Set-Volume 20 "+1.44%"

Set-Price 21 "0.996"
Set-Volume 21 "+0.16%"

Set-Price 22 "17.66"
Set-Volume 22 "-1.25%"

Set-Price 23 "5.15"
Set-Volume 23 "+0.32%"

Set-Price 24 "101.87"
Set-Volume 24 "-2.34%"

Set-Volume 25 "+0.18%"

Set-Volume 26 "-1.25%"

Set-Price 27 "10.78"
Set-Volume 27 "-1.81%"

Set-Volume 28 "-0.97%"

Set-Price 29 "33.97"
Set-Volume 29 "+0.92%"

Set-Volume 30 "+4.81%"

Set-Volume 31 "+0.29%"

Set-Price 32 "12.30"
Set-Volume 32 "-2.84%"

Set-Volume 33 "+0.81%"

Set-Price 34 "63.29"
Set-Volume 34 "+0.06%"

Set-Price 35 "0.0₃0895"
Set-Volume 35 "+8.51%"

Set-Price 36 "3.945.92"
Set-Volume 36 "+5.65%"

Set-Volume 37 "+1.44%"

Set-Price 38 "525.78"
Set-Volume 38 "+5.33%"

Set-Volume 39 "+0.03%"

Set-Price 40 "36.95"
Set-Volume 40 "+0.29%"

Set-Price 41 "0.392"
Set-Volume 41 "-0.46%"

Set-Price 42 "3.53"
Set-Volume 42 "-1.09%"

Set-Volume 43 "-1.57%"

Set-Volume 44 "-0.52%"

Set-Volume 45 "+2.25%"

Set-Volume 46 "+0.87%"

Set-Price 48 "8.62"
Set-Volume 48 "-1.13%"

Set-Volume 49 "-0.21%"

Set-Price 50 "0.000255"
Set-Volume 50 "+4.03%"

Set-Volume 51 "+3.77%"
